$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 26.144619
$ws.Range("H2").Value = 78.433857
$ws.Range("I2").Value = 0.5211737020083955
$ws.Range("J2").Value = 0.5211737020083955
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.5203476666666668
$ws.Range("N2").Value = 1.561043
$ws.Range("O2").Value = 0.004105934376266647
$ws.Range("P2").Value = 0.004105934376266647
$ws.Range("Q2").Value = 13.604291492539
$ws.Range("R2").Value = 122.438623432851
$ws.Range("S2").Value = 0.002139905019082421
$ws.Range("T2").Value = 0.002139905019082421

$ws.Range("G3").Value = 26.144619
$ws.Range("H3").Value = 78.433857
$ws.Range("I3").Value = 0.5211737020083955
$ws.Range("J3").Value = 0.5211737020083955
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 105.9632263333333
$ws.Range("N3").Value = 317.889679
$ws.Range("O3").Value = 0.8361295370252257
$ws.Range("P3").Value = 0.8361295370252259
$ws.Range("Q3").Value = 2770.368180495767
$ws.Range("R3").Value = 24933.31362446191
$ws.Range("S3").Value = 0.4357687261700027
$ws.Range("T3").Value = 0.4357687261700028

$ws.Range("G4").Value = 26.144619
$ws.Range("H4").Value = 78.433857
$ws.Range("I4").Value = 0.5211737020083955
$ws.Range("J4").Value = 0.5211737020083955
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 20.24706
$ws.Range("N4").Value = 60.74118
$ws.Range("O4").Value = 0.1597645285985076
$ws.Range("P4").Value = 0.1597645285985076
$ws.Range("Q4").Value = 529.3516695701401
$ws.Range("R4").Value = 4764.16502613126
$ws.Range("S4").Value = 0.08326507081931038
$ws.Range("T4").Value = 0.0832650708193104

$ws.Range("G5").Value = 17.91585
$ws.Range("H5").Value = 53.74755
$ws.Range("I5").Value = 0.3571392594830743
$ws.Range("J5").Value = 0.3571392594830742
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5203476666666668
$ws.Range("N5").Value = 1.561043
$ws.Range("O5").Value = 0.004105934376266647
$ws.Range("P5").Value = 0.004105934376266647
$ws.Range("Q5").Value = 9.322470743850003
$ws.Range("R5").Value = 83.90223669465001
$ws.Range("S5").Value = 0.001466390362625969
$ws.Range("T5").Value = 0.001466390362625969

$ws.Range("G6").Value = 17.91585
$ws.Range("H6").Value = 53.74755
$ws.Range("I6").Value = 0.3571392594830743
$ws.Range("J6").Value = 0.3571392594830742
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 105.9632263333333
$ws.Range("N6").Value = 317.889679
$ws.Range("O6").Value = 0.8361295370252257
$ws.Range("P6").Value = 0.8361295370252259
$ws.Range("Q6").Value = 1898.42126850405
$ws.Range("R6").Value = 17085.79141653645
$ws.Range("S6").Value = 0.2986146836851148
$ws.Range("T6").Value = 0.2986146836851148

$ws.Range("G7").Value = 17.91585
$ws.Range("H7").Value = 53.74755
$ws.Range("I7").Value = 0.3571392594830743
$ws.Range("J7").Value = 0.3571392594830742
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 20.24706
$ws.Range("N7").Value = 60.74118
$ws.Range("O7").Value = 0.1597645285985076
$ws.Range("P7").Value = 0.1597645285985076
$ws.Range("Q7").Value = 362.7432899010001
$ws.Range("R7").Value = 3264.689609109
$ws.Range("S7").Value = 0.05705818543533345
$ws.Range("T7").Value = 0.05705818543533345

$ws.Range("G8").Value = 6.104416333333333
$ws.Range("H8").Value = 18.313249
$ws.Range("I8").Value = 0.1216870385085301
$ws.Range("J8").Value = 0.1216870385085301
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.5203476666666668
$ws.Range("N8").Value = 1.561043
$ws.Range("O8").Value = 0.004105934376266647
$ws.Range("P8").Value = 0.004105934376266647
$ws.Range("Q8").Value = 3.176418795411889
$ws.Range("R8").Value = 28.587769158707
$ws.Range("S8").Value = 0.0004996389945582571
$ws.Range("T8").Value = 0.0004996389945582572

$ws.Range("G9").Value = 6.104416333333333
$ws.Range("H9").Value = 18.313249
$ws.Range("I9").Value = 0.1216870385085301
$ws.Range("J9").Value = 0.1216870385085301
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 105.9632263333333
$ws.Range("N9").Value = 317.889679
$ws.Range("O9").Value = 0.8361295370252257
$ws.Range("P9").Value = 0.8361295370252259
$ws.Range("Q9").Value = 646.8436495618968
$ws.Range("R9").Value = 5821.59284605707
$ws.Range("S9").Value = 0.1017461271701081
$ws.Range("T9").Value = 0.1017461271701081

$ws.Range("G10").Value = 6.104416333333333
$ws.Range("H10").Value = 18.313249
$ws.Range("I10").Value = 0.1216870385085301
$ws.Range("J10").Value = 0.1216870385085301
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 20.24706
$ws.Range("N10").Value = 60.74118
$ws.Range("O10").Value = 0.1597645285985076
$ws.Range("P10").Value = 0.1597645285985076
$ws.Range("Q10").Value = 123.59648376598
$ws.Range("R10").Value = 1112.36835389382
$ws.Range("S10").Value = 0.01944127234386376
$ws.Range("T10").Value = 0.01944127234386376
